$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "asis" sheet view: re-freeze/scroll the pane to the top and move
#    the selection, without leaving that sheet active afterwards.
# ------------------------------------------------------------------
$wsAsis = $wb.Worksheets.Item("asis")
$wsAsis.Activate()
$excel.ActiveWindow.ScrollRow = 2
$wsAsis.Range("B6").Select()

# ------------------------------------------------------------------
# 2) "todo" sheet: rebuild the whole task list with the new layout
#    (column A now carries an "X"/"x" marker, headers moved to column
#    B, sub-items moved to column C) and append the new todo items.
# ------------------------------------------------------------------
$wsTodo = $wb.Worksheets.Item("todo")
$wsTodo.Activate()
$wsTodo.Cells.Clear()

$wsTodo.Range("A2").Value = "X"
$wsTodo.Range("B2").Value = "graphics.plot_one_type(model, model.vars[t], emp_priors, t) <- data.plot_asr(t)"
$wsTodo.Range("B2").Font.Bold = $true

$wsTodo.Range("A3").Value = "x"
$wsTodo.Range("C3").Value = "delete data.plot_asr"

$wsTodo.Range("A4").Value = "x"
$wsTodo.Range("C4").Value = "plot_one_type(model, type, with_data, with_ui, axis)"

$wsTodo.Range("A5").Value = "x"
$wsTodo.Range("C5").Value = "delete delta displayed on plot"

$wsTodo.Range("A7").Value = "X"
$wsTodo.Range("B7").Value = "graphics.plot_one_ppc(model.vars[t],t)"
$wsTodo.Range("B7").Font.Bold = $true

$wsTodo.Range("A8").Value = "x"
$wsTodo.Range("C8").Value = "remove legend"

$wsTodo.Range("A9").Value = "x"
$wsTodo.Range("C9").Value = "plot_one_ppc(model, t)"

$wsTodo.Range("B11").Value = "graphics.plot_fit(model, model.vars, emp_priors, posterior) <- graphics.plot_cur_params(model.vars)"
$wsTodo.Range("B11").Font.Bold = $true

$wsTodo.Range("C12").Value = "delete plot_cur_params"

$wsTodo.Range("C13").Value = "axis labels"

$wsTodo.Range("C14").Value = "plot_fit(model, [type(s)], emp_priors, with_data, with_ui, axes, fig_size)"

$wsTodo.Range("A16").Value = "X"
$wsTodo.Range("B16").Value = "graphics.plot_data_bars(df, style, color)"
$wsTodo.Range("B16").Font.Bold = $true

$wsTodo.Range("A17").Value = "x"
$wsTodo.Range("C17").Value = "graphics.plot_data_bars(df, style, color, label)"

# the "delete ..." block (rows 21-25) was filled in before the
# "plot_convergence_diag" header (row 19) got its own text, so the
# shared-string table gets the same append order here.
$wsTodo.Range("A21").Value = "X"
$wsTodo.Range("A22").Value = "x"
$wsTodo.Range("C22").Value = "delete data.plot_effects"
$wsTodo.Range("B21").Value = "delete plotting in other functions"
$wsTodo.Range("B21").Font.Bold = $true

$wsTodo.Range("A23").Value = "x"
$wsTodo.Range("C23").Value = "delete data.plot_asr"

$wsTodo.Range("A24").Value = "x"
$wsTodo.Range("C24").Value = "delete data.vars.plot_trace"

$wsTodo.Range("A25").Value = "x"
$wsTodo.Range("C25").Value = "delete data.vars.plot_acorr"

$wsTodo.Range("B19").Value = "graphics.plot_convergence_diag(vars) -> graphics.plot_acorr(vars)"
$wsTodo.Range("B19").Font.Bold = $true

$wsTodo.Range("A21").Select()
